$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2 (idade) and F2 (preferencias_musicais) values
$ws.Range("D2").Value = "18/11/2004"
$ws.Range("F2").Value = "Pop, Alternativa"

# Delete row 3 entirely (removes the second user record)
$ws.Rows("3:3").Delete()
